$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new values to column A (rows 3 and 4), which Excel will
# store as shared strings appended after the existing "β" entry.
$ws.Range("A3").Value = "Côte"
$ws.Range("A4").Value = "欄外裏面に記載"
